$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1
$ws.Range("B1").Value = "Antarctic mass change (Gigatonnes)"

# Update existing data values B2:B20 (years 2002-2020)
$newValues = @(
    -6.8285714285714283,
    -115.77363636363636,
    -263.09500000000003,
    -229.1883333333333,
    -130.84416666666667,
    -317.20083333333332,
    -587.89,
    -563.11166666666657,
    -841.61916666666673,
    -940.11222222222227,
    -1074.8372727272729,
    -1285.1311111111111,
    -1495.9399999999998,
    -1847.6000000000001,
    -1769.5422222222221,
    -1803.23,
    -2147.9679999999998,
    -2269.3016666666667,
    -2542.3516666666669
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

# Add new row 21 for year 2021
$ws.Cells.Item(21, 1).Value = 2021
$ws.Cells.Item(21, 2).Value = -2913.3050000000003

# Update the selection to E17 as shown in the diff
$ws.Range("E17").Select()
